$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.42766425382235
$ws.Cells.Item(2, 3).Value = 8.427414583316333
$ws.Cells.Item(2, 4).Value = 5.963108808551331
$ws.Cells.Item(2, 5).Value = 11.98782532598147
$ws.Cells.Item(2, 7).Value = 23.17167564713594
$ws.Cells.Item(2, 8).Value = 12.5445949370599
$ws.Cells.Item(2, 13).Value = 14.03703616595114
$ws.Cells.Item(2, 14).Value = 16.41249102773353
$ws.Cells.Item(2, 15).Value = 18.39355203554458

$ws.Cells.Item(3, 2).Value = 10.86516943622082
$ws.Cells.Item(3, 3).Value = 8.088960860294682
$ws.Cells.Item(3, 4).Value = 5.840153089886527
$ws.Cells.Item(3, 5).Value = 11.87239369642501
$ws.Cells.Item(3, 7).Value = 22.973279432288
$ws.Cells.Item(3, 8).Value = 12.57489782729059
$ws.Cells.Item(3, 13).Value = 13.76052017596205
$ws.Cells.Item(3, 14).Value = 16.46694205877631
$ws.Cells.Item(3, 15).Value = 18.40119252298524

$ws.Cells.Item(4, 2).Value = 10.50551810267777
$ws.Cells.Item(4, 3).Value = 7.872365063969617
$ws.Cells.Item(4, 4).Value = 5.765116096129297
$ws.Cells.Item(4, 5).Value = 11.80559634225844
$ws.Cells.Item(4, 7).Value = 22.86104530854675
$ws.Cells.Item(4, 8).Value = 12.59622615308106
$ws.Cells.Item(4, 13).Value = 13.59111784521103
$ws.Cells.Item(4, 14).Value = 16.50217205575349
$ws.Cells.Item(4, 15).Value = 18.41145443260397

$ws.Cells.Item(5, 2).Value = 10.35554677301717
$ws.Cells.Item(5, 3).Value = 7.781974630373615
$ws.Cells.Item(5, 4).Value = 5.734701339224297
$ws.Cells.Item(5, 5).Value = 11.77943057541989
$ws.Cells.Item(5, 7).Value = 22.81777068059904
$ws.Cells.Item(5, 8).Value = 12.60560059545602
$ws.Cells.Item(5, 13).Value = 13.52227907798822
$ws.Cells.Item(5, 14).Value = 16.51698143150424
$ws.Cells.Item(5, 15).Value = 18.41703412805609

$ws.Cells.Item(6, 2).Value = 10.33044399469641
$ws.Cells.Item(6, 3).Value = 7.766839612411873
$ws.Cells.Item(6, 4).Value = 5.729662282604667
$ws.Cells.Item(6, 5).Value = 11.7751502730156
$ws.Cells.Item(6, 7).Value = 22.81073504578432
$ws.Cells.Item(6, 8).Value = 12.60719841701979
$ws.Cells.Item(6, 13).Value = 13.51086308418349
$ws.Cells.Item(6, 14).Value = 16.51946789977179
$ws.Cells.Item(6, 15).Value = 18.41804495779121

$ws.Cells.Item(7, 2).Value = 10.50350908683822
$ws.Cells.Item(7, 3).Value = 7.871154516866065
$ws.Cells.Item(7, 4).Value = 5.764705186236624
$ws.Cells.Item(7, 5).Value = 11.8052391545039
$ws.Cells.Item(7, 7).Value = 22.86045166013716
$ws.Cells.Item(7, 8).Value = 12.59634981689599
$ws.Cells.Item(7, 13).Value = 13.59018854441525
$ws.Cells.Item(7, 14).Value = 16.50236994516463
$ws.Cells.Item(7, 15).Value = 18.41152402691167

$ws.Cells.Item(8, 2).Value = 11.23677823051851
$ws.Cells.Item(8, 3).Value = 8.312585152075789
$ws.Cells.Item(8, 4).Value = 5.920647232685595
$ws.Cells.Item(8, 5).Value = 11.94719640382752
$ws.Cells.Item(8, 7).Value = 23.10131311679699
$ws.Cells.Item(8, 8).Value = 12.55447747587737
$ws.Cells.Item(8, 13).Value = 13.94167329722321
$ws.Cells.Item(8, 14).Value = 16.4308934295704
$ws.Cells.Item(8, 15).Value = 18.39502883630185

$ws.Cells.Item(9, 2).Value = 12.55492158649543
$ws.Cells.Item(9, 3).Value = 9.105432751202807
$ws.Cells.Item(9, 4).Value = 6.228013885831764
$ws.Cells.Item(9, 5).Value = 12.25652570150065
$ws.Cells.Item(9, 7).Value = 23.64707837271448
$ws.Cells.Item(9, 8).Value = 12.4940308171819
$ws.Cells.Item(9, 13).Value = 14.62971455930249
$ws.Cells.Item(9, 14).Value = 16.30493831284607
$ws.Cells.Item(9, 15).Value = 18.40696945174458

$ws.Cells.Item(10, 2).Value = 13.44326418808323
$ws.Cells.Item(10, 3).Value = 9.640082214175804
$ws.Cells.Item(10, 4).Value = 6.452111200377285
$ws.Cells.Item(10, 5).Value = 12.50062401855998
$ws.Cells.Item(10, 7).Value = 24.08897227343813
$ws.Cells.Item(10, 8).Value = 12.46291109377723
$ws.Cells.Item(10, 13).Value = 15.12900829359884
$ws.Cells.Item(10, 14).Value = 16.22099640872513
$ws.Cells.Item(10, 15).Value = 18.44279837728413

$ws.Cells.Item(11, 2).Value = 13.82890593996421
$ws.Cells.Item(11, 3).Value = 9.87235003059169
$ws.Cells.Item(11, 4).Value = 6.553153299465421
$ws.Cells.Item(11, 5).Value = 12.61485847405787
$ws.Cells.Item(11, 7).Value = 24.29797674937069
$ws.Cells.Item(11, 8).Value = 12.45165493392587
$ws.Cells.Item(11, 13).Value = 15.35371684535628
$ws.Cells.Item(11, 14).Value = 16.18466225931169
$ws.Cells.Item(11, 15).Value = 18.46496350845503

$ws.Cells.Item(12, 2).Value = 13.97221207416575
$ws.Cells.Item(12, 3).Value = 9.958691625023866
$ws.Cells.Item(12, 4).Value = 6.591244532042413
$ws.Cells.Item(12, 5).Value = 12.65853443172638
$ws.Cells.Item(12, 7).Value = 24.37818866842472
$ws.Cells.Item(12, 8).Value = 12.44781057666444
$ws.Cells.Item(12, 13).Value = 15.43837528491133
$ws.Cells.Item(12, 14).Value = 16.17116872037533
$ws.Cells.Item(12, 15).Value = 18.47419783402552

$ws.Cells.Item(13, 2).Value = 13.941470872565
$ws.Cells.Item(13, 3).Value = 9.940168721603285
$ws.Cells.Item(13, 4).Value = 6.583049155729971
$ws.Cells.Item(13, 5).Value = 12.64911011295517
$ws.Cells.Item(13, 7).Value = 24.36086749010228
$ws.Cells.Item(13, 8).Value = 12.44861991751698
$ws.Cells.Item(13, 13).Value = 15.42016317745345
$ws.Cells.Item(13, 14).Value = 16.17406300381817
$ws.Cells.Item(13, 15).Value = 18.47217171346902

$ws.Cells.Item(14, 2).Value = 13.84075086178405
$ws.Cells.Item(14, 3).Value = 9.879485938079226
$ws.Cells.Item(14, 4).Value = 6.556290749458848
$ws.Cells.Item(14, 5).Value = 12.61844358584601
$ws.Cells.Item(14, 7).Value = 24.30455493483103
$ws.Cells.Item(14, 8).Value = 12.45133027147529
$ws.Cells.Item(14, 13).Value = 15.36069091237269
$ws.Cells.Item(14, 14).Value = 16.18354682304584
$ws.Cells.Item(14, 15).Value = 18.46570638903228

$ws.Cells.Item(15, 2).Value = 13.77869976351531
$ws.Cells.Item(15, 3).Value = 9.842104842673573
$ws.Cells.Item(15, 4).Value = 6.539876960730664
$ws.Cells.Item(15, 5).Value = 12.59971263215293
$ws.Cells.Item(15, 7).Value = 24.27019828911317
$ws.Cells.Item(15, 8).Value = 12.4530449194684
$ws.Cells.Item(15, 13).Value = 15.32420348523838
$ws.Cells.Item(15, 14).Value = 16.18939047708218
$ws.Cells.Item(15, 15).Value = 18.46185559851677

$ws.Cells.Item(16, 2).Value = 13.41768346239673
$ws.Cells.Item(16, 3).Value = 9.624679207338806
$ws.Cells.Item(16, 4).Value = 6.44548611541665
$ws.Cells.Item(16, 5).Value = 12.49321926445521
$ws.Cells.Item(16, 7).Value = 24.07546707952098
$ws.Cells.Item(16, 8).Value = 12.46370516518568
$ws.Cells.Item(16, 13).Value = 15.1142668971927
$ws.Cells.Item(16, 14).Value = 16.22340811592531
$ws.Cells.Item(16, 15).Value = 18.44146767383633

$ws.Cells.Item(17, 2).Value = 13.19142545165643
$ws.Cells.Item(17, 3).Value = 9.488462080410228
$ws.Cells.Item(17, 4).Value = 6.387319229467114
$ws.Cells.Item(17, 5).Value = 12.42867767413763
$ws.Cells.Item(17, 7).Value = 23.9579934930295
$ws.Cells.Item(17, 8).Value = 12.47098851944452
$ws.Cells.Item(17, 13).Value = 14.98479478170975
$ws.Cells.Item(17, 14).Value = 16.24475044317103
$ws.Cells.Item(17, 15).Value = 18.43046128218158

$ws.Cells.Item(18, 2).Value = 13.05955274089238
$ws.Cells.Item(18, 3).Value = 9.40908512641993
$ws.Cells.Item(18, 4).Value = 6.353780720812042
$ws.Cells.Item(18, 5).Value = 12.39185811827453
$ws.Cells.Item(18, 7).Value = 23.89118242176513
$ws.Cells.Item(18, 8).Value = 12.47545065509169
$ws.Cells.Item(18, 13).Value = 14.91010262362235
$ws.Cells.Item(18, 14).Value = 16.25720029942785
$ws.Cells.Item(18, 15).Value = 18.42468309874456

$ws.Cells.Item(19, 2).Value = 13.01460732383294
$ws.Cells.Item(19, 3).Value = 9.382034105851925
$ws.Cells.Item(19, 4).Value = 6.34241230604458
$ws.Cells.Item(19, 5).Value = 12.37944494620439
$ws.Cells.Item(19, 7).Value = 23.86869374370952
$ws.Cells.Item(19, 8).Value = 12.47700830163291
$ws.Cells.Item(19, 13).Value = 14.88477748601569
$ws.Cells.Item(19, 14).Value = 16.26144557087321
$ws.Cells.Item(19, 15).Value = 18.42282164531164

$ws.Cells.Item(20, 2).Value = 13.21569114285936
$ws.Cells.Item(20, 3).Value = 9.503069370799153
$ws.Cells.Item(20, 4).Value = 6.393520037125848
$ws.Cells.Item(20, 5).Value = 12.43551716793509
$ws.Cells.Item(20, 7).Value = 23.97042103327062
$ws.Cells.Item(20, 8).Value = 12.47018493735427
$ws.Cells.Item(20, 13).Value = 14.99860104629887
$ws.Cells.Item(20, 14).Value = 16.24246048115189
$ws.Cells.Item(20, 15).Value = 18.43157577594787

$ws.Cells.Item(21, 2).Value = 13.87040931762898
$ws.Cells.Item(21, 3).Value = 9.897354004810607
$ws.Cells.Item(21, 4).Value = 6.564155308665957
$ws.Cells.Item(21, 5).Value = 12.62744008016078
$ws.Cells.Item(21, 7).Value = 24.32106701631034
$ws.Cells.Item(21, 8).Value = 12.45052282002974
$ws.Cells.Item(21, 13).Value = 15.37817178845241
$ws.Cells.Item(21, 14).Value = 16.18075399858007
$ws.Cells.Item(21, 15).Value = 18.46758261749272

$ws.Cells.Item(22, 2).Value = 14.2823799059461
$ws.Cells.Item(22, 3).Value = 10.14562497094
$ws.Cells.Item(22, 4).Value = 6.674659770180429
$ws.Cells.Item(22, 5).Value = 12.75528880822487
$ws.Cells.Item(22, 7).Value = 24.55641126340637
$ws.Cells.Item(22, 8).Value = 12.44010990367533
$ws.Cells.Item(22, 13).Value = 15.62367847997707
$ws.Cells.Item(22, 14).Value = 16.14197192286946
$ws.Cells.Item(22, 15).Value = 18.49601438944837

$ws.Cells.Item(23, 2).Value = 14.06398114496208
$ws.Cells.Item(23, 3).Value = 10.01399107535219
$ws.Cells.Item(23, 4).Value = 6.615787506783169
$ws.Cells.Item(23, 5).Value = 12.68684626048118
$ws.Cells.Item(23, 7).Value = 24.43026596931571
$ws.Cells.Item(23, 8).Value = 12.44544414690263
$ws.Cells.Item(23, 13).Value = 15.49290856519704
$ws.Cells.Item(23, 14).Value = 16.16252939018471
$ws.Cells.Item(23, 15).Value = 18.48039271382816

$ws.Cells.Item(24, 2).Value = 13.20472620606558
$ws.Cells.Item(24, 3).Value = 9.496468725697484
$ws.Cells.Item(24, 4).Value = 6.390716952779756
$ws.Cells.Item(24, 5).Value = 12.43242413797642
$ws.Cells.Item(24, 7).Value = 23.96480027448105
$ws.Cells.Item(24, 8).Value = 12.47054738076608
$ws.Cells.Item(24, 13).Value = 14.99236002774585
$ws.Cells.Item(24, 14).Value = 16.24349521261638
$ws.Cells.Item(24, 15).Value = 18.43107020131342

$ws.Cells.Item(25, 2).Value = 12.21196513715645
$ws.Cells.Item(25, 3).Value = 8.899132852032704
$ws.Cells.Item(25, 4).Value = 6.144978736980099
$ws.Cells.Item(25, 5).Value = 12.16973325586557
$ws.Cells.Item(25, 7).Value = 23.49195148184513
$ws.Cells.Item(25, 8).Value = 12.50805477561766
$ws.Cells.Item(25, 13).Value = 14.44428956515366
$ws.Cells.Item(25, 14).Value = 16.33749808402536
$ws.Cells.Item(25, 15).Value = 18.3989878827321
